$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 4.5
$ws.Range("J8").Value = 3
$ws.Range("L8").Value = 9
$ws.Range("N8").Value = -287
$ws.Range("H17").Value = 560.8125
$ws.Range("J17").Value = 553.129
$ws.Range("L17").Value = 1659.387
$ws.Range("N17").Value = -1995.387
$ws.Range("H116").Value = 10713
$ws.Range("J116").Value = 13197.5
$ws.Range("L116").Value = 13197.5
$ws.Range("N116").Value = -20081.5
$ws.Range("H132").Value = 6993.2
$ws.Range("I132").Value = 7994.8237
$ws.Range("J132").Value = 1317.3334
$ws.Range("K132").Value = 23984.4711
$ws.Range("L132").Value = 3952.0002
$ws.Range("M132").Value = -21454.4711
$ws.Range("N132").Value = -9012.0002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2118064
$ws.Range("I32").Value = 1030570.1
$ws.Range("K32").Value = 1030570.1
$ws.Range("M32").Value = -1030283.1
$ws.Range("H45").Value = 31837.572
$ws.Range("I45").Value = 43831.7
$ws.Range("J45").Value = 1852.25
$ws.Range("K45").Value = 43831.7
$ws.Range("L45").Value = 1852.25
$ws.Range("M45").Value = -43454.7
$ws.Range("N45").Value = -2606.25
$ws.Range("H61").Value = 2489.6191
$ws.Range("I61").Value = 2237.9443
$ws.Range("J61").Value = 3999.6667
$ws.Range("K61").Value = 2237.9443
$ws.Range("L61").Value = 3999.6667
$ws.Range("M61").Value = -2025.9443
$ws.Range("N61").Value = -4423.6667
$ws.Range("H97").Value = 559.619
$ws.Range("I97").Value = 475.5
$ws.Range("J97").Value = 828.8
$ws.Range("K97").Value = 475.5
$ws.Range("L97").Value = 828.8
$ws.Range("M97").Value = 20.5
$ws.Range("N97").Value = -1820.8
$ws.Range("H122").Value = 3043.625
$ws.Range("I122").Value = 3452.6924
$ws.Range("K122").Value = 10358.0772
$ws.Range("M122").Value = -7908.0772
$ws.Range("H132").Value = 2381.3914
$ws.Range("I132").Value = 1851.5333
$ws.Range("K132").Value = 5554.5999
$ws.Range("M132").Value = -3024.5999
$ws.Range("H136").Value = 2489.6191
$ws.Range("I136").Value = 2237.9443
$ws.Range("J136").Value = 3999.6667
$ws.Range("K136").Value = 6713.8329
$ws.Range("L136").Value = 11999.0001
$ws.Range("M136").Value = -4163.8329
$ws.Range("N136").Value = -17099.0001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3434.5334
$ws.Range("I86").Value = 3492.5454
$ws.Range("J86").Value = 3275
$ws.Range("K86").Value = 3492.5454
$ws.Range("L86").Value = 3275
$ws.Range("M86").Value = -2369.5454
$ws.Range("N86").Value = -5521
$ws.Range("H89").Value = 3434.5334
$ws.Range("I89").Value = 3492.5454
$ws.Range("J89").Value = 3275
$ws.Range("K89").Value = 17462.727
$ws.Range("L89").Value = 16375
$ws.Range("M89").Value = -11846.727
$ws.Range("N89").Value = -27607
$ws.Range("H94").Value = 133340264
$ws.Range("I94").Value = 142864420
$ws.Range("K94").Value = 142864420
$ws.Range("M94").Value = -142863969
$ws.Range("H99").Value = 2568.1538
$ws.Range("I99").Value = 1674.125
$ws.Range("J99").Value = 3998.6
$ws.Range("K99").Value = 1674.125
$ws.Range("L99").Value = 3998.6
$ws.Range("M99").Value = -176.125
$ws.Range("N99").Value = -6994.6
$ws.Range("H134").Value = 2554.9048
$ws.Range("I134").Value = 1796.3636
$ws.Range("K134").Value = 5389.0908
$ws.Range("M134").Value = -2854.0908
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1136.8
$ws.Range("I16").Value = 1046.25
$ws.Range("K16").Value = 1046.25
$ws.Range("M16").Value = -759.25
$ws.Range("H22").Value = 492.5
$ws.Range("I22").Value = 391
$ws.Range("K22").Value = 391
$ws.Range("M22").Value = -41
$ws.Range("H31").Value = 6947376
$ws.Range("J31").Value = 31253454
$ws.Range("L31").Value = 31253454
$ws.Range("N31").Value = -31254044
$ws.Range("H34").Value = 6947376
$ws.Range("J34").Value = 31253454
$ws.Range("L34").Value = 31253454
$ws.Range("N34").Value = -31253858
$ws.Range("H99").Value = 4155.357
$ws.Range("I99").Value = 3562.25
$ws.Range("K99").Value = 3562.25
$ws.Range("M99").Value = -2064.25
$ws.Range("H105").Value = 2449.2
$ws.Range("I105").Value = 2449.2
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2449.2
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H107").Value = 3125974.5
$ws.Range("I107").Value = 5000680.5
$ws.Range("K107").Value = 5000680.5
$ws.Range("M107").Value = -4998760.5
$ws.Range("H113").Value = 1136.8
$ws.Range("I113").Value = 1046.25
$ws.Range("K113").Value = 1046.25
$ws.Range("M113").Value = 1123.75
$ws.Range("H122").Value = 2354.3076
$ws.Range("I122").Value = 2093.3333
$ws.Range("J122").Value = 3450.4
$ws.Range("K122").Value = 6279.999899999999
$ws.Range("L122").Value = 10351.2
$ws.Range("M122").Value = -3829.999899999999
$ws.Range("N122").Value = -15251.2
$ws.Range("H126").Value = 4155.357
$ws.Range("I126").Value = 3562.25
$ws.Range("K126").Value = 10686.75
$ws.Range("M126").Value = -8216.75
$ws.Range("H134").Value = 4587.593
$ws.Range("J134").Value = 3699.25
$ws.Range("L134").Value = 11097.75
$ws.Range("N134").Value = -16167.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 42500
$ws.Range("J23").Value = 42500
$ws.Range("L23").Value = 127500
$ws.Range("N23").Value = -127970
$ws.Range("H68").Value = 9097162
$ws.Range("I68").Value = 1073.75
$ws.Range("K68").Value = 3221.25
$ws.Range("M68").Value = -2410.25
$ws.Range("H71").Value = 9097162
$ws.Range("I71").Value = 1073.75
$ws.Range("K71").Value = 9663.75
$ws.Range("M71").Value = -5607.75
$ws.Range("H122").Value = 1161.6
$ws.Range("I122").Value = 603
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 5427
$ws.Range("L122").Value = 17995.5
$ws.Range("M122").Value = -2977
$ws.Range("N122").Value = -22895.5
$ws.Range("H140").Value = 6605.8
$ws.Range("I140").Value = 6605.8
$ws.Range("K140").Value = 19817.4
$ws.Range("M140").Value = -14637.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 166670820
$ws.Range("I80").Value = 500001760
$ws.Range("J80").Value = 5337.5
$ws.Range("K80").Value = 500001760
$ws.Range("L80").Value = 5337.5
$ws.Range("M80").Value = -500000762
$ws.Range("N80").Value = -7333.5
$ws.Range("H83").Value = 166670820
$ws.Range("I83").Value = 500001760
$ws.Range("J83").Value = 5337.5
$ws.Range("K83").Value = 2500008800
$ws.Range("L83").Value = 26687.5
$ws.Range("M83").Value = -2500003808
$ws.Range("N83").Value = -36671.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3140.7144
$ws.Range("I7").Value = 2996.25
$ws.Range("J7").Value = 3333.3333
$ws.Range("K7").Value = 2996.25
$ws.Range("L7").Value = 3333.3333
$ws.Range("M7").Value = -2884.25
$ws.Range("N7").Value = -3557.3333
$ws.Range("H93").Value = 2804.3635
$ws.Range("I93").Value = 2612.75
$ws.Range("J93").Value = 3315.3333
$ws.Range("K93").Value = 2612.75
$ws.Range("L93").Value = 3315.3333
$ws.Range("M93").Value = -1364.75
$ws.Range("N93").Value = -5811.3333
$ws.Range("H122").Value = 6999.25
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 6999.25
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("N122").Value = -25897.75
$ws.Range("H126").Value = 3140.7144
$ws.Range("I126").Value = 2996.25
$ws.Range("J126").Value = 3333.3333
$ws.Range("K126").Value = 8988.75
$ws.Range("L126").Value = 9999.999899999999
$ws.Range("M126").Value = -6518.75
$ws.Range("N126").Value = -14939.9999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3751.2
$ws.Range("I81").Value = 4407.5557
$ws.Range("J81").Value = 3214.182
$ws.Range("K81").Value = 8815.1114
$ws.Range("L81").Value = 6428.364
$ws.Range("M81").Value = -7754.1114
$ws.Range("N81").Value = -8550.364
$ws.Range("H84").Value = 3751.2
$ws.Range("I84").Value = 4407.5557
$ws.Range("J84").Value = 3214.182
$ws.Range("K84").Value = 48861.25
$ws.Range("L84").Value = 34456
$ws.Range("M84").Value = -43557.25
$ws.Range("N84").Value = -45064
$ws.Range("H101").Value = 32940.832
$ws.Range("J101").Value = 32940.832
$ws.Range("L101").Value = 32940.832
$ws.Range("N101").Value = -39430.832
$ws.Range("H122").Value = 13891010
$ws.Range("I122").Value = 2265.3333
$ws.Range("K122").Value = 6795.999899999999
$ws.Range("M122").Value = -4345.999899999999
$ws.Range("H126").Value = 12500.182
$ws.Range("I126").Value = 15781.75
$ws.Range("K126").Value = 47345.25
$ws.Range("M126").Value = -44875.25
$ws.Range("H136").Value = 3186.2
$ws.Range("I136").Value = 3270.9285
$ws.Range("K136").Value = 9812.7855
$ws.Range("M136").Value = -7262.7855
